$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.366.17'
$ws.Range("E2").Value = '  -5.92%  '
$ws.Range("D3").Value = '2.448.59'
$ws.Range("E3").Value = '  -8.83%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '534.65'
$ws.Range("E5").Value = '  -2.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.63'
$ws.Range("E6").Value = '  -6.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0991'
$ws.Range("E9").Value = '  -6.23%  '
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.38'
$ws.Range("E11").Value = '  +5.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.351'
$ws.Range("E12").Value = '  -5.03%  '
$ws.Range("D13").Value = '2.879.22'
$ws.Range("E13").Value = '  -8.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.17'
$ws.Range("E14").Value = '  -7.58%  '
$ws.Range("D15").Value = '59.326.40'
$ws.Range("E15").Value = '  -5.76%  '
$ws.Range("E16").Value = '  -6.28%  '
$ws.Range("D17").Value = '2.492.78'
$ws.Range("E17").Value = '  -7.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.12'
$ws.Range("E18").Value = '  -6.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.37'
$ws.Range("E19").Value = '  -4.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '325.00'
$ws.Range("E20").Value = '  -5.50%  '
$ws.Range("E21").Value = '  -3.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.69'
$ws.Range("E22").Value = '  -9.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.468'
$ws.Range("E23").Value = '  -7.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.38'
$ws.Range("E24").Value = '  -5.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.162'
$ws.Range("E25").Value = '  -3.80%  '
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.76'
$ws.Range("E27").Value = '  -4.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.92'
$ws.Range("E28").Value = '  -1.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.30'
$ws.Range("E29").Value = '  -2.73%  '
$ws.Range("E30").Value = '  -5.61%  '
$ws.Range("D31").Value = '0.0₃0773'
$ws.Range("E31").Value = '  -10.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.997'
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '156.86'
$ws.Range("E33").Value = '  -6.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.55'
$ws.Range("E34").Value = '  -5.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.31'
$ws.Range("E35").Value = '  -6.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.35'
$ws.Range("E36").Value = '  -5.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.73'
$ws.Range("E37").Value = '  -2.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.78'
$ws.Range("E38").Value = '  -6.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '312.86'
$ws.Range("E39").Value = '  -7.87%  '
$ws.Range("E40").Value = '  -8.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.77'
$ws.Range("E41").Value = '  -3.91%  '
$ws.Range("E42").Value = '  -6.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.994'
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.71'
$ws.Range("E44").Value = '  -3.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.583'
$ws.Range("E45").Value = '  -5.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0939'
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0523'
$ws.Range("E47").Value = '  -6.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.09'
$ws.Range("E48").Value = '  -8.37%  '
$ws.Range("E49").Value = '  -4.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.50'
$ws.Range("E50").Value = '  -9.25%  '
$ws.Range("D51").Value = '1.988.35'
$ws.Range("E51").Value = '  -5.03%  '
